$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the six fonts whose names move around once the list is sorted ---
# (sizes stay put; only the font name used at that size/row changes)
$ws.Range("A3").Value = "Arial"
$ws.Range("A3").Font.Name = "Arial"

$ws.Range("A4").Value = "Arial Black"
$ws.Range("A4").Font.Name = "Arial Black"

$ws.Range("A6").Value = "Arial Rounded MT Bold"
$ws.Range("A6").Font.Name = "Arial Rounded MT Bold"

$ws.Range("A7").Value = "Bahnschrift"
$ws.Range("A7").Font.Name = "Bahnschrift"

$ws.Range("A8").Value = "Barlow Condensed"
$ws.Range("A8").Font.Name = "Barlow Condensed"

$ws.Range("A10").Value = "Barlow Condensed ExtraBold"
$ws.Range("A10").Font.Name = "Barlow Condensed ExtraBold"

# --- Append ten more example rows, each with its own font name/size ---
$newFonts = @(
    @{Row=11; Name="Barlow Condensed ExtraLight"; Size=30},
    @{Row=12; Name="Barlow Condensed Light";      Size=32},
    @{Row=13; Name="Barlow Condensed Medium";     Size=34},
    @{Row=14; Name="Barlow Condensed SemiBold";   Size=36},
    @{Row=15; Name="Barlow Condensed Thin";       Size=38},
    @{Row=16; Name="Baskerville Old Face";        Size=40},
    @{Row=17; Name="Bauhaus 93";                  Size=42},
    @{Row=18; Name="Bell MT";                     Size=44},
    @{Row=19; Name="Berlin Sans FB";               Size=46},
    @{Row=20; Name="Berlin Sans FB Demi";          Size=48}
)

foreach ($item in $newFonts) {
    $cell = $ws.Cells.Item($item.Row, 1)
    $cell.Value = $item.Name
    $cell.Font.Size = $item.Size
    $cell.Font.Name = $item.Name
}
